$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 347
$ws.Range("I33").Value = 269.5
$ws.Range("K33").Value = 269.5
$ws.Range("M33").Value = -40.5
$ws.Range("H40").Value = 6770
$ws.Range("I40").Value = 4621.1113
$ws.Range("J40").Value = 9993.333000000001
$ws.Range("K40").Value = 4621.1113
$ws.Range("L40").Value = 9993.333000000001
$ws.Range("M40").Value = -4446.1113
$ws.Range("N40").Value = -10343.333
$ws.Range("H61").Value = 3999
$ws.Range("I61").Value = 3999
$ws.Range("K61").Value = 11997
$ws.Range("M61").Value = -11825
$ws.Range("H92").Value = 641.75
$ws.Range("I92").Value = 647.2143
$ws.Range("K92").Value = 647.2143
$ws.Range("M92").Value = 600.7857
$ws.Range("H99").Value = 244.22223
$ws.Range("I99").Value = 199.76471
$ws.Range("K99").Value = 599.29413
$ws.Range("M99").Value = 898.70587
$ws.Range("H113").Value = 3801.2
$ws.Range("I113").Value = 2699
$ws.Range("K113").Value = 2699
$ws.Range("M113").Value = 555
$ws.Range("H138").Value = 2963.8513
$ws.Range("I138").Value = 2383.3845
$ws.Range("J138").Value = 3278.2708
$ws.Range("K138").Value = 7150.1535
$ws.Range("L138").Value = 9834.812399999999
$ws.Range("M138").Value = -2010.1535
$ws.Range("N138").Value = -20114.8124

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3570.5881
$ws.Range("I122").Value = 3668.75
$ws.Range("K122").Value = 11006.25
$ws.Range("M122").Value = -8556.25
$ws.Range("H132").Value = 2718.125
$ws.Range("J132").Value = 3666.3333
$ws.Range("L132").Value = 10998.9999
$ws.Range("N132").Value = -16058.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 153847120
$ws.Range("I94").Value = 181818980
$ws.Range("K94").Value = 181818980
$ws.Range("M94").Value = -181818529

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1618.4546
$ws.Range("I16").Value = 709.75
$ws.Range("J16").Value = 2137.7144
$ws.Range("K16").Value = 709.75
$ws.Range("L16").Value = 2137.7144
$ws.Range("M16").Value = -422.75
$ws.Range("N16").Value = -2711.7144
$ws.Range("H58").Value = 2191.2593
$ws.Range("I58").Value = 1308.55
$ws.Range("K58").Value = 1308.55
$ws.Range("M58").Value = -1105.55
$ws.Range("H62").Value = 12502357
$ws.Range("I62").Value = 16668827
$ws.Range("K62").Value = 16668827
$ws.Range("M62").Value = -16668203
$ws.Range("H65").Value = 12502357
$ws.Range("I65").Value = 16668827
$ws.Range("K65").Value = 83344135
$ws.Range("M65").Value = -83341015
$ws.Range("H107").Value = 2778728.8
$ws.Range("I107").Value = 5555943
$ws.Range("K107").Value = 5555943
$ws.Range("M107").Value = -5554023
$ws.Range("H113").Value = 1618.4546
$ws.Range("I113").Value = 709.75
$ws.Range("J113").Value = 2137.7144
$ws.Range("K113").Value = 709.75
$ws.Range("L113").Value = 2137.7144
$ws.Range("M113").Value = 1460.25
$ws.Range("N113").Value = -6477.7144
$ws.Range("H132").Value = 15878636
$ws.Range("I132").Value = 3694.3845
$ws.Range("K132").Value = 11083.1535
$ws.Range("M132").Value = -8553.1535
$ws.Range("H134").Value = 2220.16
$ws.Range("I134").Value = 2009.1333
$ws.Range("K134").Value = 6027.3999
$ws.Range("M134").Value = -3492.3999
$ws.Range("H136").Value = 2191.2593
$ws.Range("I136").Value = 1308.55
$ws.Range("K136").Value = 3925.65
$ws.Range("M136").Value = -1375.65

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 211.85185
$ws.Range("J2").Value = 333
$ws.Range("L2").Value = 1998
$ws.Range("N2").Value = -2224
$ws.Range("H55").Value = 4685.636
$ws.Range("J55").Value = 5605.3125
$ws.Range("L55").Value = 16815.9375
$ws.Range("N55").Value = -17169.9375
$ws.Range("H93").Value = 8666.666999999999
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("H103").Value = 1884.3334
$ws.Range("J103").Value = 5000
$ws.Range("L103").Value = 15000
$ws.Range("N103").Value = -16758

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 533.3333
$ws.Range("I2").Value = 300
$ws.Range("K2").Value = 300
$ws.Range("M2").Value = -187
$ws.Range("H11").Value = 2165437
$ws.Range("I11").Value = 1152
$ws.Range("J11").Value = 2420058.8
$ws.Range("K11").Value = 1152
$ws.Range("L11").Value = 2420058.8
$ws.Range("M11").Value = -1013
$ws.Range("N11").Value = -2420336.8
$ws.Range("H122").Value = 7697460.5
$ws.Range("I122").Value = 12825430
$ws.Range("J122").Value = 5506
$ws.Range("K122").Value = 38476290
$ws.Range("L122").Value = 16518
$ws.Range("M122").Value = -38473840
$ws.Range("N122").Value = -21418
$ws.Range("H132").Value = 2582.6667
$ws.Range("I132").Value = 2421.2856
$ws.Range("J132").Value = 2959.2222
$ws.Range("K132").Value = 7263.8568
$ws.Range("L132").Value = 8877.6666
$ws.Range("M132").Value = -4733.8568
$ws.Range("N132").Value = -13937.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8240.375
$ws.Range("J40").Value = 8751
$ws.Range("L40").Value = 8751
$ws.Range("N40").Value = -9023
$ws.Range("H132").Value = 5392.7144
$ws.Range("I132").Value = 1775.1818
$ws.Range("J132").Value = 11514.692
$ws.Range("K132").Value = 5325.5454
$ws.Range("L132").Value = 34544.076
$ws.Range("M132").Value = -2795.5454
$ws.Range("N132").Value = -39604.076

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 41348
$ws.Range("I42").Value = 24044
$ws.Range("J42").Value = 50000
$ws.Range("K42").Value = 24044
$ws.Range("L42").Value = 50000
$ws.Range("M42").Value = -23666
$ws.Range("N42").Value = -50756
$ws.Range("H96").Value = 2484.8
$ws.Range("I96").Value = 1364
$ws.Range("K96").Value = 1364
$ws.Range("M96").Value = 9
$ws.Range("H100").Value = 66668170
$ws.Range("I100").Value = 1753.1
$ws.Range("K100").Value = 3506.2
$ws.Range("M100").Value = -2965.2
$ws.Range("H113").Value = 1602.8125
$ws.Range("I113").Value = 1444.7
$ws.Range("J113").Value = 1866.3334
$ws.Range("K113").Value = 4334.1
$ws.Range("L113").Value = 5599.0002
$ws.Range("M113").Value = -2164.1
$ws.Range("N113").Value = -9939.0002
$ws.Range("H125").Value = 148999
$ws.Range("J125").Value = 148999
$ws.Range("L125").Value = 148999
$ws.Range("N125").Value = -158839
$ws.Range("H132").Value = 2620.465
$ws.Range("I132").Value = 2409.8975
$ws.Range("J132").Value = 4673.5
$ws.Range("K132").Value = 7229.6925
$ws.Range("L132").Value = 14020.5
$ws.Range("M132").Value = -4699.6925
$ws.Range("N132").Value = -19080.5
